$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data at the end of the table (row 19)
# Copy formatting (incl. date number format) from the row above to avoid
# introducing a duplicate/custom number format.
$ws.Cells.Item(18, 1).Copy() | Out-Null
$ws.Cells.Item(19, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(19, 1).Value = 45529
$ws.Cells.Item(19, 2).Value = 120000
$ws.Cells.Item(19, 3).Value = "Paul en River"
$ws.Cells.Item(19, 4).Value = 45225577
$ws.Cells.Item(19, 5).Value = "Lucas Pablo"
$ws.Cells.Item(19, 6).Value = "Antiñolo"
$ws.Cells.Item(19, 7).Value = "2x Tipo: Platea Preferencial"

# Resize the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:G19"))
